$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

# Remove the obsolete question row (row 5: "Is dit een goede nieuwe vraag??").
# This shifts the remaining rows (old 6 -> 5, old 7 -> 6) up by one, matching
# the new dimension A1:L6.
$ws.Rows.Item(5).Delete()

# The question that is now in row 5 (originally row 6, "4 Wat is er aan de
# gang?") gets a newly generated image_url.
$ws.Range("L5").Value = "https://raw.githubusercontent.com/onomatorHanze/didactic-octo-spork/main/data/images/DC_new_1763126867.png"
